$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "species"

# --- Row 2: CHEMBL230 ---
$ws.Range("B2").Value = "Cyclooxygenase-2"
$ws.Range("C2").Value = "homo sapiens"

# --- Row 3: CHEMBL25 ---
$ws.Range("B3").Value = "aspirin"
$ws.Range("C3").Value = "na"

# --- Row 4: CHEMBL4523582 ---
$ws.Range("B4").Value = "Replicase polyprotein 1ab"
$ws.Range("C4").Value = "Severe acute respiratory syndrome coronavirus 2 "

# --- Row 5: CHEMBL3616356 ---
$ws.Range("B5").Value = "D-aspartate oxidase"
$ws.Range("C5").Value = "Rattus norvegicus "

# --- Row 6: new row, CHEMBL1871 ---
$ws.Range("A6").Value = "CHEMBL1871"
$ws.Range("B6").Value = "androgen receptor"
# Reuse the exact same shared string as C2 ("homo sapiens") by copy/paste
$ws.Range("C2").Copy($ws.Range("C6"))

# --- Apply word-wrap styling to the cells that need it ---
$ws.Range("B2").WrapText = $true
$ws.Range("B3").WrapText = $true
$ws.Range("B4").WrapText = $true
$ws.Range("C4").WrapText = $true
$ws.Range("B5").WrapText = $true
$ws.Range("C5").WrapText = $true
$ws.Range("A6").WrapText = $true
$ws.Range("B6").WrapText = $true

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 23.85
$ws.Rows.Item(4).RowHeight = 46.25
$ws.Rows.Item(5).RowHeight = 23.85
$ws.Rows.Item(6).RowHeight = 23.85

# --- Column B width ---
$ws.Columns.Item(2).ColumnWidth = 13.8

# --- Selection moves to A7 after editing ---
[void]$ws.Range("A7").Select()
